# Quarterly update for TSC financials: add new quarter columns (Dec-2018, Sep-2018)
# as D:E, shifting the prior 8 quarters right by two columns, and correct two
# Capital Expenditures data points (Mar-2017, Jun-2017) that were re-stated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column D; this shifts the existing
# quarterly data (old D:K) to F:M.
$ws.Columns("D:E").Insert()

# New columns D and E need the same per-row formatting (date format on the
# "Period Ending" rows, numeric format elsewhere) as the data that now sits
# in column F. Copy that formatting down in one shot for each new column.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# (row, D-value, E-value) for every row that carries data in the new quarter
# columns. $null means leave the cell blank (it already is, post-insert).
$rowData = @(
    @(7, 43465, 43373),
    @(8, 58200, 52400),
    @(9, "NA", "NA"),
    @(10, "NA", "NA"),
    @(12, "NA", "NA"),
    @(13, 0, 0),
    @(14, 0, 0),
    @(15, -500, -500),
    @(17, 28000, 23400),
    @(18, 30200, 29000),
    @(20, -14800, -12900),
    @(21, 16300, 17000),
    @(22, 0, 0),
    @(23, 15400, 16100),
    @(24, 600, 1800),
    @(25, 0, 0),
    @(26, 14800, 14300),
    @(27, 14100, 13600),
    @(28, 0, 0),
    @(29, 300, "NA"),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 14800, 12900),
    @(33, 14400, 13600),
    @(34, 0, 0),
    @(35, 14400, 13600),
    @(38, 43465, 43373),
    @(41, 500, 400),
    @(42, 214300, 203000),
    @(43, 0, 0),
    @(44, 0, 0),
    @(45, 0, 0),
    @(46, 0, 0),
    @(47, 0, 0),
    @(48, 5100, 5000),
    @(49, 67900, 68400),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, "NA", "NA"),
    @(53, 0, 0),
    @(54, 6035700, 5573300),
    @(57, 5200, 4100),
    @(58, 0, 0),
    @(59, 90000, 76800),
    @(60, 0, 0),
    @(61, 0, 0),
    @(62, 3500, 4800),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 5556300, 5105700),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 38500, 38500),
    @(71, 0, 0),
    @(72, 164000, 149600),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 440900, 429200),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 14400, 13600),
    @(83, 900, 900),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 27100, 23100),
    @(91, -500, -500),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -457800, -309500),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, 434100, 308600),
    @(101, 0, 0),
    @(102, 3500, 22200)

)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($r, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($r, 5).Value = $eVal
    }
}

# Two Capital Expenditures (row 91) values were also restated for the
# Mar-2017 and Jun-2017 quarters, now in columns I and J after the shift.
$ws.Cells.Item(91, 9).Value = -200
$ws.Cells.Item(91, 10).Value = -300

Write-Output "done"
